# Financial model update: insert a new reporting-period column (D) ahead of
# the existing yearly data on the WGO worksheet, preserving formats, and
# populate it with the latest fiscal year (period ending 2018-08-25,
# serial date 43337) figures across the Income Statement, Balance Sheet
# and Cash Flow Statement sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WGO")

# Insert a new column before column D; this shifts the existing D:K data
# right to E:L (column widths/dimension follow automatically).
$ws.Columns("D").Insert()

# The freshly inserted column D has no number formatting yet (Excel seeds
# it from the column to the left). Copy formatting only from column E so
# the new cells pick up the correct date / number styles.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Set numeric values for the new column D (fiscal year ending 2018-08-25, serial 43337)
$newValues = @{
    7 = 43337
    8 = 2016800
    9 = 1717000
    10 = 299800
    13 = 0
    15 = 9300
    17 = 1856400
    18 = 160400
    20 = 500
    21 = 180100
    22 = 18200
    23 = 142600
    24 = 40300
    25 = 0
    26 = 102400
    27 = 102400
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -500
    33 = 102400
    34 = 0
    35 = 102400
    38 = 43337
    41 = 2300
    42 = 0
    43 = 164600
    44 = 195100
    45 = 9900
    46 = 371900
    47 = 28300
    48 = 103200
    49 = 538100
    50 = 0
    51 = 0
    52 = 10300
    53 = 0
    54 = 1051800
    57 = 81000
    58 = 0
    59 = 123100
    60 = 204200
    61 = 291400
    62 = 21700
    63 = 0
    64 = 0
    65 = 0
    66 = 517400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 768800
    73 = 0
    74 = 0
    75 = 0
    76 = 534400
    77 = 0
    80 = 43337
    81 = 102400
    83 = 19200
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 83300
    91 = -28700
    92 = 0
    93 = 0
    94 = -111800
    96 = -12700
    97 = 0
    98 = 0
    99 = 0
    100 = -5200
    101 = 0
    102 = -33600
}
foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}

# Rows where the new period has no data available ("NA")
$naRows = @(12, 14)
foreach ($row in $naRows) {
    $ws.Range("D$row").Value = "NA"
}
